$d = $word.ActiveDocument
$sec = $d.Sections(1)

# The document has two BTec (orange) logos living in the header stories and
# two Pearson logos living in the footer stories. Each is an inline picture
# whose drawing "name" (wp:docPr / pic:cNvPr @name) needs to be swapped:
#   BTec logos:    image2.jpg -> image1.jpg
#   Pearson logos: image1.png -> image2.png
#
# InlineShape has no settable Name in the Word object model, so each
# picture is briefly converted to a floating Shape (which does expose
# .Name), renamed, then converted back to an inline picture in place.

function Rename-InlineLogo($range, $newName) {
    $inline = $range.InlineShapes(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape()
}

Rename-InlineLogo $sec.Headers(1).Range "image1.jpg"
Rename-InlineLogo $sec.Headers(2).Range "image1.jpg"
Rename-InlineLogo $sec.Footers(1).Range "image2.png"
Rename-InlineLogo $sec.Footers(2).Range "image2.png"
